$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: existing experiment FE9251 (exp 4) now gets a description logged
$ws.Range("A6").Value = 4
$ws.Range("K6").Value = "Agrego la tendencia de 12 meses en dataset de train"

# New "Descripcion" column header
$ws.Range("K1").Value = "Descripcion"

# Row 7: new experiment FE9252 (exp 5)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "FE9252"
$ws.Range("C7").Value = "z925_FE_historia-v1"
$ws.Range("K7").Value = "Activo fn CanaritosAsesinos"

# Size the new column to fit its (long) text, like the other columns in the sheet
$ws.Columns.Item(11).AutoFit()

# Leave the selection where the last edit happened
$ws.Range("D7").Select()
